$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Opt Portfolio" (C) and "Opt Portfolio with View" (D) columns
# with the newly computed weights for rows 2-8.
$ws.Range("C2:D2").Value = 0.2718729031447604
$ws.Range("C3:D3").Value = 0.2619866744857898
$ws.Range("C4:D4").Value = 0
$ws.Range("C5:D5").Value = [double]"2.846030702774449e-19"
$ws.Range("C6:D6").Value = [double]"3.469446951953614e-18"
$ws.Range("C7:D7").Value = 0.2407114716596564
$ws.Range("C8:D8").Value = 0.2254289507097933
